$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3000
$ws.Range("D2").Value = 4500
$ws.Range("D3").Value = 4000

$ws.Range("D4").Select()
